$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 19; existing rows 19-144 shift down to 20-145.
$ws.Rows.Item(19).Insert()

# Fill the new row 19 with data (copy constant columns, set new varying values).
$ws.Cells.Item(19, 1).Value = 3
$ws.Cells.Item(19, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(19, 3).Value = "Coquimbo"
$ws.Cells.Item(19, 4).Value = 44613
$ws.Cells.Item(19, 5).Value = 5
$ws.Cells.Item(19, 6).Value = 100112052
$ws.Cells.Item(19, 7).Value = "Albahaca"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 125
$ws.Cells.Item(19, 11).Value = 5000
$ws.Cells.Item(19, 12).Value = 5500
$ws.Cells.Item(19, 13).Value = 5240
$ws.Cells.Item(19, 14).Value = "`$/docena de matas"
$ws.Cells.Item(19, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(19, 16).Value = 873
$ws.Cells.Item(19, 17).Value = 6
$ws.Cells.Item(19, 18).Value = "Hortaliza"
